# Update Name of Algo
# Apply updated numeric results to the RandomForest result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 15.53220000000001
$ws.Range("E4").Value = 16.3413
$ws.Range("C6").Value = -11.79400000000001
$ws.Range("B7").Value = 5.8345
$ws.Range("D7").Value = -7.637899999999998
$ws.Range("A9").Value = -21.89330000000002
$ws.Range("D10").Value = -7.9889
$ws.Range("B12").Value = 5.383499999999995
$ws.Range("D13").Value = -8.164800000000005
$ws.Range("E13").Value = 16.1964
$ws.Range("B14").Value = 5.917800000000002
$ws.Range("C15").Value = -14.37869999999999
$ws.Range("D16").Value = -8.801500000000006
$ws.Range("A18").Value = -22.28860000000002
$ws.Range("A20").Value = -20.40369999999998
$ws.Range("D20").Value = -7.115799999999996
$ws.Range("D24").Value = -7.3787
$ws.Range("B26").Value = 4.015500000000005
$ws.Range("A27").Value = -21.87789999999999
$ws.Range("B27").Value = 5.757500000000001
$ws.Range("E27").Value = 16.32019999999999
$ws.Range("B29").Value = 5.119499999999995
$ws.Range("E29").Value = 17.22710000000001
$ws.Range("C33").Value = -11.37369999999999
$ws.Range("A35").Value = -19.52039999999999
$ws.Range("C35").Value = -12.63650000000001
$ws.Range("E35").Value = 16.2627
$ws.Range("B37").Value = 9.323900000000007
$ws.Range("B38").Value = 4.875799999999999
$ws.Range("C38").Value = -12.4789
$ws.Range("D39").Value = -7.302400000000002
$ws.Range("E40").Value = 17.05030000000001
$ws.Range("C43").Value = -13.8559
$ws.Range("C44").Value = -13.2462
$ws.Range("C47").Value = -12.0179
$ws.Range("D47").Value = -7.219399999999998
$ws.Range("D48").Value = -7.117699999999996
$ws.Range("B51").Value = 6.260100000000006
$ws.Range("C51").Value = -11.9309
$ws.Range("B52").Value = 5.509099999999998
$ws.Range("D52").Value = -7.296799999999993
$ws.Range("B55").Value = 4.833399999999997
$ws.Range("D56").Value = -7.666899999999997
$ws.Range("C57").Value = -14.42959999999999
$ws.Range("E57").Value = 16.6216
$ws.Range("C63").Value = -11.838
$ws.Range("A69").Value = -21.6461
$ws.Range("B69").Value = 5.327799999999996
$ws.Range("B70").Value = 6.019700000000008
$ws.Range("C70").Value = -11.4202
$ws.Range("A76").Value = -19.82569999999998
$ws.Range("A78").Value = -20.01099999999998
$ws.Range("B81").Value = 5.545600000000003
$ws.Range("A82").Value = -21.94240000000001
$ws.Range("A83").Value = -20.75439999999997
$ws.Range("B83").Value = 7.267399999999996
$ws.Range("D84").Value = -9.246999999999998
$ws.Range("E85").Value = 16.2031
$ws.Range("C88").Value = -12.0828
$ws.Range("A93").Value = -20.76789999999998
$ws.Range("C99").Value = -12.453
$ws.Range("D100").Value = -8.459200000000003
$ws.Range("D101").Value = -8.000699999999997
$ws.Range("B102").Value = 8.290700000000003
